$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Add Devices Loop A")
$ws2 = $wb.Worksheets.Item("Other Devices Loop A")
$ws3 = $wb.Worksheets.Item("Sheet1")

# --- Sheet "Add Devices Loop A": add I1:I3, mirroring E1:E3's formatting ---
$ws1.Range("E1:E3").Copy()
$ws1.Range("I1:I3").PasteSpecial(-4122)
$ws1.Range("I1").Value = "DC Unit Loading Details Name"
$ws1.Range("I2").Value = "Current (DC Units)"
$ws1.Range("I3").Value = "Current (worst case)"

# --- Sheet "Other Devices Loop A": add I1:I3, mirroring E1:E3's formatting ---
$ws2.Range("E1:E3").Copy()
$ws2.Range("I1:I3").PasteSpecial(-4122)
$ws2.Range("I1").Value = "DC Unit Loading Details Name"
$ws2.Range("I2").Value = "Current (DC Units)"
$ws2.Range("I3").Value = "Current (worst case)"

# --- Update selections on every sheet; select the final sheet last so it
#     ends up the active tab (matches "Other Devices Loop A" staying active).
$ws1.Range("I1:I3").Select()
$ws3.Range("I1").Select()
$ws2.Range("I1:I3").Select()
